# Auto-generated Excel COM-interop script
# Applies the IESO report update: refreshed CreatedAt timestamp and updated
# numeric values in columns U-Z across the data rows (GitHub Actions report refresh).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-08-04T17:07:44"
$ws.Range("U4").Value = 116.3
$ws.Range("V4").Value = 99
$ws.Range("W4").Value = 71.78
$ws.Range("X4").Value = 60
$ws.Range("Y4").Value = 38.93
$ws.Range("Z4").Value = 34.35
$ws.Range("U6").Value = -6.63
$ws.Range("X6").Value = -2.85
$ws.Range("Y6").Value = -1.17
$ws.Range("Z6").Value = -0.65
$ws.Range("X8").Value = -9.6
$ws.Range("U9").Value = 116.41
$ws.Range("V9").Value = 98.81
$ws.Range("W9").Value = 72.68000000000001
$ws.Range("X9").Value = 62.21
$ws.Range("Y9").Value = 40.29
$ws.Range("Z9").Value = 35.97
$ws.Range("U11").Value = -6.52
$ws.Range("V11").Value = -5.83
$ws.Range("W11").Value = -2.83
$ws.Range("X11").Value = -0.65
$ws.Range("Y11").Value = 0.2
$ws.Range("Z11").Value = 0.97
$ws.Range("X13").Value = -9.6
$ws.Range("U14").Value = 116.41
$ws.Range("V14").Value = 98.81
$ws.Range("W14").Value = 72.75
$ws.Range("X14").Value = 62.28
$ws.Range("Y14").Value = 40.33
$ws.Range("Z14").Value = 35.97
$ws.Range("U16").Value = -6.52
$ws.Range("V16").Value = -5.83
$ws.Range("W16").Value = -2.76
$ws.Range("X16").Value = -0.57
$ws.Range("Y16").Value = 0.24
$ws.Range("Z16").Value = 0.97
$ws.Range("X18").Value = -9.6
$ws.Range("U19").Value = 115.87
$ws.Range("V19").Value = 98.53
$ws.Range("W19").Value = 71.44
$ws.Range("X19").Value = 59.73
$ws.Range("Y19").Value = 38.74
$ws.Range("U21").Value = -7.07
$ws.Range("V21").Value = -6.11
$ws.Range("X21").Value = -3.12
$ws.Range("Y21").Value = -1.36
$ws.Range("X23").Value = -9.6
$ws.Range("U24").Value = 115.87
$ws.Range("V24").Value = 98.53
$ws.Range("W24").Value = 71.44
$ws.Range("X24").Value = 59.73
$ws.Range("Y24").Value = 38.74
$ws.Range("U26").Value = -7.07
$ws.Range("V26").Value = -6.11
$ws.Range("X26").Value = -3.12
$ws.Range("Y26").Value = -1.36
$ws.Range("X28").Value = -9.6
$ws.Range("U29").Value = 115.21
$ws.Range("V29").Value = 97.89
$ws.Range("W29").Value = 70.91
$ws.Range("X29").Value = 59.21
$ws.Range("Y29").Value = 38.4
$ws.Range("Z29").Value = 34.18
$ws.Range("U31").Value = -7.72
$ws.Range("V31").Value = -6.75
$ws.Range("W31").Value = -4.61
$ws.Range("X31").Value = -3.65
$ws.Range("Y31").Value = -1.69
$ws.Range("Z31").Value = -0.82
$ws.Range("X33").Value = -9.6
$ws.Range("U34").Value = 117.53
$ws.Range("V34").Value = 99.56999999999999
$ws.Range("W34").Value = 73.81999999999999
$ws.Range("X34").Value = 64.03
$ws.Range("Y34").Value = 41.29
$ws.Range("Z34").Value = 37.04
$ws.Range("U36").Value = -5.41
$ws.Range("V36").Value = -5.08
$ws.Range("W36").Value = -1.7
$ws.Range("X36").Value = 1.18
$ws.Range("Y36").Value = 1.2
$ws.Range("Z36").Value = 2.04
$ws.Range("X38").Value = -9.6
$ws.Range("U39").Value = 116.3
$ws.Range("V39").Value = 99
$ws.Range("W39").Value = 71.78
$ws.Range("X39").Value = 60
$ws.Range("Y39").Value = 38.93
$ws.Range("Z39").Value = 34.35
$ws.Range("U41").Value = -6.63
$ws.Range("X41").Value = -2.85
$ws.Range("Y41").Value = -1.17
$ws.Range("Z41").Value = -0.65
$ws.Range("X43").Value = -9.6
$ws.Range("U44").Value = 119.01
$ws.Range("V44").Value = 102.09
$ws.Range("W44").Value = 74.11
$ws.Range("X44").Value = 61.64
$ws.Range("Y44").Value = 39.97
$ws.Range("U46").Value = -3.93
$ws.Range("V46").Value = -2.55
$ws.Range("X46").Value = -1.21
$ws.Range("Y46").Value = -0.12
$ws.Range("X48").Value = -9.6
$ws.Range("U49").Value = 121.24
$ws.Range("V49").Value = 102.79
$ws.Range("X49").Value = 62.85
$ws.Range("Y49").Value = 40.42
$ws.Range("U51").Value = -1.7
$ws.Range("V51").Value = -1.85
$ws.Range("W51").Value = -0.6
$ws.Range("X51").Value = 0
$ws.Range("Y51").Value = 0.32
$ws.Range("X53").Value = -9.6
$ws.Range("U54").Value = 115.54
$ws.Range("V54").Value = 99.47
$ws.Range("W54").Value = 72.68000000000001
$ws.Range("X54").Value = 61.57
$ws.Range("Y54").Value = 40.25
$ws.Range("U56").Value = -7.39
$ws.Range("V56").Value = -5.17
$ws.Range("X56").Value = -1.28
$ws.Range("Y56").Value = 0.16
$ws.Range("X58").Value = -9.6
$ws.Range("U59").Value = 122.81
$ws.Range("V59").Value = 105.59
$ws.Range("W59").Value = 76.75
$ws.Range("X59").Value = 64.03
$ws.Range("Y59").Value = 41.25
$ws.Range("V61").Value = 0.95
$ws.Range("X61").Value = 1.18
$ws.Range("Y61").Value = 1.15
$ws.Range("X63").Value = -9.6
$ws.Range("U64").Value = 125.06
$ws.Range("V64").Value = 107.66
$ws.Range("W64").Value = 78.26000000000001
$ws.Range("X64").Value = 65.48
$ws.Range("Y64").Value = 41.89
$ws.Range("U66").Value = 2.13
$ws.Range("V66").Value = 3.01
$ws.Range("X66").Value = 2.63
$ws.Range("Y66").Value = 1.8
$ws.Range("X68").Value = -9.6
$ws.Range("U69").Value = 123.55
$ws.Range("V69").Value = 106.89
$ws.Range("W69").Value = 78.09
$ws.Range("X69").Value = 65.33
$ws.Range("Y69").Value = 42.25
$ws.Range("Z69").Value = 37
$ws.Range("U71").Value = 0.62
$ws.Range("V71").Value = 2.24
$ws.Range("W71").Value = 2.58
$ws.Range("X71").Value = 2.47
$ws.Range("Y71").Value = 2.15
$ws.Range("Z71").Value = 2
$ws.Range("X73").Value = -9.6
$ws.Range("U74").Value = 122.32
$ws.Range("V74").Value = 105.27
$ws.Range("W74").Value = 76.43000000000001
$ws.Range("X74").Value = 63.73
$ws.Range("Y74").Value = 41.04
$ws.Range("U76").Value = -0.61
$ws.Range("V76").Value = 0.63
$ws.Range("X76").Value = 0.88
$ws.Range("Y76").Value = 0.9399999999999999
$ws.Range("X78").Value = -9.6
$ws.Range("U79").Value = 123.04
$ws.Range("V79").Value = 105.95
$ws.Range("W79").Value = 76.95999999999999
$ws.Range("X79").Value = 64.37
$ws.Range("Y79").Value = 41.32
$ws.Range("Z79").Value = 36.1
$ws.Range("U81").Value = 0.11
$ws.Range("V81").Value = 1.31
$ws.Range("W81").Value = 1.44
$ws.Range("X81").Value = 1.52
$ws.Range("Y81").Value = 1.22
$ws.Range("Z81").Value = 1.1
$ws.Range("X83").Value = -9.6
$ws.Range("U84").Value = 112.37
$ws.Range("V84").Value = 97.98
$ws.Range("W84").Value = 73.39
$ws.Range("X84").Value = 62.35
$ws.Range("Y84").Value = 40.42
$ws.Range("Z84").Value = 35.35
$ws.Range("U86").Value = -10.56
$ws.Range("V86").Value = -6.66
$ws.Range("W86").Value = -2.13
$ws.Range("X86").Value = -0.5
$ws.Range("Y86").Value = 0.32
$ws.Range("Z86").Value = 0.35
$ws.Range("X88").Value = -9.6
$ws.Range("U89").Value = 115.21
$ws.Range("V89").Value = 97.89
$ws.Range("W89").Value = 70.91
$ws.Range("X89").Value = 59.21
$ws.Range("Y89").Value = 38.4
$ws.Range("Z89").Value = 34.18
$ws.Range("U91").Value = -7.72
$ws.Range("V91").Value = -6.75
$ws.Range("W91").Value = -4.61
$ws.Range("X91").Value = -3.65
$ws.Range("Y91").Value = -1.69
$ws.Range("Z91").Value = -0.82
$ws.Range("X93").Value = -9.6
